$d = $word.ActiveDocument

# 1) The paragraph "The following improvements to our webpage were performed
#    this week." used to be split across two runs: a leading-space-only run
#    and the sentence run. Re-run a Find/Replace over the whole sentence
#    (including the leading space) so Word collapses it back into a single
#    run, matching the target markup.
$d.Content.Find.Execute(
    " The following improvements to our webpage were performed this week.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " The following improvements to our webpage were performed this week.",
    2) | Out-Null

# 2) Append the new weekly-update paragraphs at the very end of the story,
#    each separated by a blank paragraph, exactly like the rest of the
#    document. The last paragraph reproduces the grammar-checker markup
#    (<w:proofErr>) that wraps "All of" in the target document, so we build
#    it from a raw WordprocessingML fragment and insert it with InsertXML.
$endPos = $d.Content.End
$insertionRange = $d.Range($endPos, $endPos)

$newBodyXml = '<w:p/>' + `
    '<w:p><w:r><w:t>I worked on the application''s coding with my team, and we were able to construct a home page with all the titles and background colors.</w:t></w:r></w:p>' + `
    '<w:p/>' + `
    '<w:p><w:r><w:t xml:space="preserve"> We also made a login page for users to register and access the page with the insurance details, where we could look at the insurance details and then choose the payment option.</w:t></w:r></w:p>' + `
    '<w:p/>' + `
    '<w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>All of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the coding has been completed up to this point; the remaining work and coding will be completed later.</w:t></w:r></w:p>'

$xmlChunk = '<?xml version="1.0" standalone="yes"?>' + `
    '<?mso-application progid="Word.Document"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $newBodyXml + '</w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertionRange.InsertXML($xmlChunk)
